$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.315.23"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.159.13"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.90"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.27"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  +6.07%  "
$ws.Range("E8").Value = "  +3.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.156.40"
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.746"
$ws.Range("E11").Value = "  +5.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.202"
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.32"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.59"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.181.18"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.157.40"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.76"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.24"
$ws.Range("E20").Value = "  +11.54%  "
$ws.Range("E21").Value = "  +12.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.62"
$ws.Range("E22").Value = "  +6.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000204"
$ws.Range("E23").Value = "  -4.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.21"
$ws.Range("E24").Value = "  +6.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.02"
$ws.Range("E25").Value = "  +8.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.09"
$ws.Range("E26").Value = "  +2.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.07"
$ws.Range("E27").Value = "  +3.87%  "
$ws.Range("E28").Value = "  +2.71%  "
$ws.Range("E30").Value = "  +47.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.236"
$ws.Range("E31").Value = "  +18.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.171"
$ws.Range("E32").Value = "  +10.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.42"
$ws.Range("E33").Value = "  +4.92%  "
$ws.Range("E34").Value = "  +15.29%  "
$ws.Range("E35").Value = "  -10.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.70"
$ws.Range("E36").Value = "  +10.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.52"
$ws.Range("E37").Value = "  +3.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "513.47"
$ws.Range("E38").Value = "  +5.16%  "
$ws.Range("E39").Value = "  +5.13%  "
$ws.Range("E40").Value = "  +9.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.90"
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.451"
$ws.Range("E42").Value = "  +13.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.49"
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.14"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.721"
$ws.Range("E46").Value = "  +7.62%  "
$ws.Range("E47").Value = "  +6.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "158.67"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("E49").Value = "  +7.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.52"
$ws.Range("E50").Value = "  +5.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.05"
$ws.Range("E51").Value = "  -0.21%  "
